$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187598586082458
$ws.Range("B1").Value = 2.1838538646698
$ws.Range("C1").Value = 6.379278659820557
$ws.Range("D1").Value = 2.303624868392944
$ws.Range("E1").Value = 1.19311511516571
